$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 105; $r -le 120; $r++) {
    $ws.Range("B$r").Value = "yes"
    $ws.Range("E$r").Value = "AIP_exchanges.feature"
    $ws.Range("F$r").Value = "AIPExchanges.java"
}

$ws.Range("F121").Select()
